$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "purpose" column (E2:E22) from "S.GISH" to the new value "fullRNASEQ"
$ws.Range("E2:E22").Value = "fullRNASEQ"

# Update the active selection to match the post-edit state (scrolled down, single cell selected)
$ws.Range("E23").Select()
